# Update latest output (run 171)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1): recomputed Cost ($) / Unit Cost ($/ML) ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 336.60947775
$schedule.Range("F2").Value = 7.420843865740741

$schedule.Range("E3").Value = 734.1133020000001
$schedule.Range("F3").Value = 27.74426689342404

$schedule.Range("E4").Value = 89.53332374999999
$schedule.Range("F4").Value = 2.631784942680776

# --- Sheet "Detailed" (sheet2): updated Price values / Type labels ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B33").Value = -10
$detailed.Range("B34").Value = -9.710129999999999

$detailed.Range("C35").Value = "historical"

$detailed.Range("B36").Value = 3.68083
$detailed.Range("C36").Value = "historical"

$detailed.Range("B37").Value = 7.44246
$detailed.Range("B38").Value = 11.67904
$detailed.Range("B39").Value = 19.83023
$detailed.Range("B40").Value = 27.5939
$detailed.Range("B41").Value = 57.31
$detailed.Range("B43").Value = 53.50268
$detailed.Range("B44").Value = 43.06138
$detailed.Range("B45").Value = 53.16109
$detailed.Range("B46").Value = 50.06385
$detailed.Range("B51").Value = 57.06
$detailed.Range("B54").Value = 48.46707
$detailed.Range("B55").Value = 48.13188
$detailed.Range("B56").Value = 49.74196
$detailed.Range("B57").Value = 56.98
$detailed.Range("B58").Value = 56.98
$detailed.Range("B60").Value = 57.06
$detailed.Range("B61").Value = 59.92295
$detailed.Range("B62").Value = 59.45974
$detailed.Range("B64").Value = 26.82262
$detailed.Range("B65").Value = 24.58885
$detailed.Range("B69").Value = 0.61
$detailed.Range("B70").Value = 0.51
$detailed.Range("B71").Value = 0.7
$detailed.Range("B72").Value = 1.22723
$detailed.Range("B73").Value = 0.7
$detailed.Range("B74").Value = 4.88872
$detailed.Range("B75").Value = 20.18987
$detailed.Range("B76").Value = 22.07
$detailed.Range("B77").Value = 0.7
$detailed.Range("B78").Value = 0.02674
$detailed.Range("B79").Value = -2.54301
$detailed.Range("B81").Value = -5.17224
$detailed.Range("B82").Value = -4.30529
$detailed.Range("B83").Value = -6.31978
$detailed.Range("B85").Value = 0.009509999999999999
$detailed.Range("B86").Value = 12.30462
$detailed.Range("B87").Value = 25.51607
$detailed.Range("B92").Value = 55.14118
$detailed.Range("B93").Value = 56.58094
$detailed.Range("B97").Value = 56.98
